$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A5").Value = "Otro"
$ws.Range("A5").Select()
